# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Update the "Rule" label of the R40 row (row 11) on the Rules sheet
# from "R40" to "1". The value must remain text (not a number), so it
# is entered with a leading apostrophe the same way a user typing into
# the cell would force text storage.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("B11").Value = "'1"
